# ProdPageTestData.xlsx update
# - prodPage sheet: drop the now-unused column D, fix the mojibake'd product
#   title text, trim the trailing space off "Write a review", and rewrap
#   the title row a bit taller.
# - jumpTo sheet: new data set for the "jump to section" test cases, with a
#   new sectionHeading column.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$xlLeft = -4131
$xlCenter = -4108

# ---------------------------------------------------------------------
# Sheet 1: prodPage
# ---------------------------------------------------------------------

# Column D is empty/unused now - remove it entirely.
$ws1.Columns.Item(4).Delete()

# Fix the product title text (was double UTF-8 encoded) and taller row.
$ws1.Range("C3").Value = "BENADRYLÂ® Allergy ULTRATABSÂ® Tablets with Allergy Relief and Diphenhydramine HCI 25 mg"
$ws1.Rows.Item(3).RowHeight = 60

# Drop the trailing space in the review button label.
$ws1.Range("C6").Value = "Write a review"

[void]$ws1.Range("C10").Select()

# ---------------------------------------------------------------------
# Sheet 2: jumpTo
# ---------------------------------------------------------------------

# Header row gets a new "sectionHeading" column; expectedText -> jumpToHeading.
$ws2.Range("C1").Value = "jumpToHeading"
$ws2.Range("D1").Value = "sectionHeading"

# Insert the 7 new test-case rows plus the new column D.
$ws2.Rows.Item(2).Resize(7).Insert()
$ws2.Columns.Item(4).Insert()

$rows = @(
  @("Verify jumping to Overview",    "OVERVIEW",  "Product Overview",           30),
  @("Verify jumping to Directions",  "DIRECTIONS","Directions",                 30),
  @("Verify jumping to Ingredients", "INGREDIENTS","Ingredients",               30),
  @("Verify jumping to Used For",    "USED FOR",  "Used For",                   30),
  @("Verify jumping to Warnings",    "WARNINGS",  "Warnings",                   30),
  @("Verify jumping to FAQs",        "FAQS",      "Frequently Asked Questions", 0),
  @("Verify jumping to Review",      "REVIEWS",   "Reviews",                    30)
)

$r = 2
foreach ($row in $rows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = "yes"
    $ws2.Cells.Item($r, 3).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = $row[2]

    $ws2.Cells.Item($r, 1).HorizontalAlignment = $xlLeft
    $ws2.Cells.Item($r, 1).VerticalAlignment = $xlCenter
    $ws2.Cells.Item($r, 1).WrapText = $true

    $ws2.Cells.Item($r, 2).HorizontalAlignment = $xlLeft
    $ws2.Cells.Item($r, 2).VerticalAlignment = $xlCenter
    $ws2.Cells.Item($r, 2).WrapText = $false

    $ws2.Cells.Item($r, 3).ClearFormats()
    $ws2.Cells.Item($r, 4).ClearFormats()

    if ($row[3] -gt 0) {
        $ws2.Rows.Item($r).RowHeight = $row[3]
    }

    $r = $r + 1
}

$ws2.Range("D1").ClearFormats()

[void]$ws2.Range("E12").Select()

Write-Host "edit complete"
